$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D text-looking numeric values to remain as text (avoid Excel auto-numeric conversion)
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue 'D2' '22.442.76'
$ws.Range("E2").Value = '  +0.46%  '

# Row 3
Set-TextValue 'D3' '1.571.40'
$ws.Range("E3").Value = '  +0.53%  '

# Row 4
$ws.Range("E4").Value = '  -0.47%  '

# Row 5
$ws.Range("E5").Value = '  -0.44%  '

# Row 6
Set-TextValue 'D6' '290.02'
$ws.Range("E6").Value = '  +0.31%  '

# Row 7
Set-TextValue 'D7' '0.3697'
$ws.Range("E7").Value = '  -1.09%  '

# Row 8
Set-TextValue 'D8' '49.86'
$ws.Range("E8").Value = '  +0.83%  '

# Row 9
Set-TextValue 'D9' '0.3388'
$ws.Range("E9").Value = '  +0.83%  '

# Row 10
Set-TextValue 'D10' '1.150'
$ws.Range("E10").Value = '  +3.14%  '

# Row 11
Set-TextValue 'D11' '0.07568'
$ws.Range("E11").Value = '  +1.44%  '

# Row 12
Set-TextValue 'D12' '1.001'
$ws.Range("E12").Value = '  -0.58%  '

# Row 13
Set-TextValue 'D13' '21.23'
$ws.Range("E13").Value = '  +2.63%  '

# Row 14
Set-TextValue 'D14' '6.027'
$ws.Range("E14").Value = '  +2.95%  '

# Row 15
Set-TextValue 'D15' '6.995'
$ws.Range("E15").Value = '  +2.16%  '

# Row 16
Set-TextValue 'D16' '1.570.61'
$ws.Range("E16").Value = '  +0.40%  '

# Row 17
Set-TextValue 'D17' '0.00001123'
$ws.Range("E17").Value = '  +1.84%  '

# Row 18
Set-TextValue 'D18' '90.58'
$ws.Range("E18").Value = '  +1.77%  '

# Row 19
Set-TextValue 'D19' '0.06776'
$ws.Range("E19").Value = '  +0.93%  '

# Row 20
$ws.Range("E20").Value = '  -0.42%  '

# Row 21
Set-TextValue 'D21' '6.368'
$ws.Range("E21").Value = '  +3.71%  '

# Row 22
Set-TextValue 'D22' '16.43'
$ws.Range("E22").Value = '  +1.32%  '

# Row 23
$ws.Range("E23").Value = '  +3.33%  '

# Row 24
Set-TextValue 'D24' '22.442.92'
$ws.Range("E24").Value = '  +0.48%  '

# Row 25
Set-TextValue 'D25' '2.370'
$ws.Range("E25").Value = '  -0.16%  '

# Row 26
Set-TextValue 'D26' '2.680'
$ws.Range("E26").Value = '  +3.44%  '

# Row 27
Set-TextValue 'D27' '20.02'
$ws.Range("E27").Value = '  +0.87%  '

# Row 28
Set-TextValue 'D28' '149.51'
$ws.Range("E28").Value = '  +1.52%  '

# Row 29
Set-TextValue 'D29' '5.049'
$ws.Range("E29").Value = '  +0.67%  '

# Row 30
Set-TextValue 'D30' '125.16'
$ws.Range("E30").Value = '  +0.60%  '

# Row 31
Set-TextValue 'D31' '1.748.56'
$ws.Range("E31").Value = '  +0.49%  '

# Row 32
Set-TextValue 'D32' '1.065'
$ws.Range("E32").Value = '  +8.98%  '

# Row 33
$ws.Range("E33").Value = '  +6.11%  '

# Row 34
Set-TextValue 'D34' '2.014'
$ws.Range("E34").Value = '  -0.24%  '

# Row 35
Set-TextValue 'D35' '9.820'
$ws.Range("E35").Value = '  +0.84%  '

# Row 36
Set-TextValue 'D36' '0.08372'
$ws.Range("E36").Value = '  -0.63%  '

# Row 37
Set-TextValue 'D37' '0.02480'
$ws.Range("E37").Value = '  +1.93%  '

# Row 38
$ws.Range("B38").Value = 'Algorand'
$ws.Range("C38").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue 'D38' '0.2301'
$ws.Range("E38").Value = '  +2.38%  '

# Row 39
$ws.Range("B39").Value = 'TrustWalletToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 'D39' '1.349'
$ws.Range("E39").Value = '  -3.82%  '

# Row 40
Set-TextValue 'D40' '0.06574'
$ws.Range("E40").Value = '  +3.24%  '

# Row 41
Set-TextValue 'D41' '5.436'
$ws.Range("E41").Value = '  +2.21%  '

# Row 42
Set-TextValue 'D42' '11.34'
$ws.Range("E42").Value = '  +4.23%  '

# Row 43
Set-TextValue 'D43' '0.6263'
$ws.Range("E43").Value = '  +1.39%  '

# Row 44
$ws.Range("B44").Value = 'Frax'
$ws.Range("C44").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue 'D44' '1.001'
$ws.Range("E44").Value = '  -0.33%  '

# Row 45
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D45' '14.07'
$ws.Range("E45").Value = '  +2.38%  '

# Row 46
Set-TextValue 'D46' '3.799'
$ws.Range("E46").Value = '  +0.41%  '

# Row 47
Set-TextValue 'D47' '0.5881'
$ws.Range("E47").Value = '  +2.52%  '

# Row 48
Set-TextValue 'D48' '2.074'
$ws.Range("E48").Value = '  +2.29%  '

# Row 49
Set-TextValue 'D49' '127.89'
$ws.Range("E49").Value = '  +3.34%  '

# Row 50
$ws.Range("E50").Value = '  +0.52%  '

# Row 51
Set-TextValue 'D51' '0.07302'
$ws.Range("E51").Value = '  +0.18%  '
